# Minor corrections to licence files
#
# The commit:
#  - Merges a few runs that were split apart only because of an
#    inline <w:proofErr> spell-check marker (e.g. "Created by " / "DMagic" /
#    ", Paul Kingtiger and Simon Hinton" -> one run), dropping the now
#    redundant <w:proofErr> markers.
#  - Changes "Italian translation provided by Carlo Rossi" to
#    "Italian translation provided by CRL42" (a differently formatted
#    second run).
#  - Moves the stray "_GoBack" bookmark from just before the "3rd Part
#    mods..." heading run down to just after "Cybutek" (i.e. where the
#    cursor was when the document was last saved).
#
# We rebuild each affected paragraph's XML exactly (preserving its
# paragraph properties) and splice it in with Range.InsertXML so the
# resulting markup matches precisely, rather than relying on
# Find/Replace's run-splitting side effects.

$d = $word.ActiveDocument

function Set-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range
    $pkg = '<?xml version="1.0" standalone="yes"?>' +
           '<?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $innerXml + '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# 1) "Created by DMagic, Paul Kingtiger and Simon Hinton" (title page byline)
$xml2 = '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Created by DMagic, Paul Kingtiger and Simon Hinton</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml 2 $xml2

# 2) "DMagic Orbital Science" heading
$xml17 = '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr>' +
    '<w:r><w:t>DMagic Orbital Science</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml 17 $xml17

# 3) "Universal Storage is created by DMagic, Paul Kingtiger and Simon Hinton (Daishi)"
$xml31 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Universal Storage is created by DMagic, Paul Kingtiger and Simon Hinton (Daishi)</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml 31 $xml31

# 4) "Italian translation provided by Carlo Rossi" -> "... CRL42"
$xml33 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Italian translation provided by </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t>CRL42</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml 33 $xml33

# 5) "AVC online created by Cybutek" -- insert the _GoBack bookmark right
#    after the "Cybutek" run (before the closing spellEnd marker).
$xml35 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">AVC online created by </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Cybutek</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
Set-ParagraphXml 35 $xml35

# 6) "Copyright 2018 DMagic, Paul Kingtiger and Simon Hinton" -- merge the
#    first two runs, leave the trailing ", Paul Kingtiger and Simon Hinton"
#    run as-is.
$xml47 = '<w:p>' +
    '<w:r><w:t>Copyright 2018 DMagic</w:t></w:r>' +
    '<w:r><w:t>, Paul Kingtiger and Simon Hinton</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml 47 $xml47

# 7) Remove the (now moved) _GoBack bookmark from before the "3rd Part
#    mods..." heading.
$xml55 = '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr>' +
    '<w:r><w:t>3</w:t></w:r>' +
    '<w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>rd</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Part mods and plug-ins packaged with Universal Storage</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml 55 $xml55
